# The three data rows 9-11 had their field values cyclically rotated:
#   new row 9  gets the values that were in row 11
#   new row 10 gets the values that were in row 9
#   new row 11 gets the values that were in row 10
# (row/column formatting stays where it is - only the field values move.)
#
# We read the current values with Value2 (plain data, no COM type-coercion
# surprises), build the rotated block in memory, and then write back only
# the cells whose value actually changes - this keeps the edit minimal and
# avoids Excel's automatic text->date/time reinterpretation on cells whose
# rotated value happens to be identical to what's already there (e.g. the
# "Startdatum"/"Slutdatum" columns, which hold the same date in all three
# rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 9
$lastRow = 11
$lastCol = 51  # column AY
$rowCount = $lastRow - $firstRow + 1

$rng = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$current = $rng.Value2   # 1-based [row, col] array: current[1..rowCount, 1..lastCol]

for ($c = 1; $c -le $lastCol; $c++) {
    for ($r = 1; $r -le $rowCount; $r++) {
        # Row r (1-based offset within the block) should end up holding what
        # used to be one row above it, wrapping from the top back to the
        # bottom (row 1 <- old row $rowCount).
        $srcR = $r - 1
        if ($srcR -lt 1) { $srcR = $rowCount }

        $oldVal = $current[$r, $c]
        $newVal = $current[$srcR, $c]

        $changed = $false
        if ($oldVal -eq $null -and $newVal -ne $null) { $changed = $true }
        elseif ($oldVal -ne $null -and $newVal -eq $null) { $changed = $true }
        elseif ($oldVal -ne $newVal) { $changed = $true }

        if ($changed) {
            $cell = $ws.Cells.Item($firstRow + $r - 1, $c)
            if ($newVal -eq $null) {
                $cell.ClearContents()
            } else {
                $cell.Value2 = $newVal
            }
        }
    }
}
